# Fix Training Data Issue
# Column BF ("Date") held malformed date strings like "6-9-2012-13"
# (an artefact of how NBA stats were originally formatted: "M-D-YYYY-YY").
# Re-write them as a proper ISO-ish "YYYY-MM-DD" text value, e.g. "2013-06-09".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "6-9-2012-13"
$newDate = "2013-06-09"

# Column BF is column 58 (B=2 ... BF=58). Data rows run 2..31.
$firstRow = 2
$lastRow  = 31
$dateCol  = 58

$dateRange = $ws.Range($ws.Cells.Item($firstRow, $dateCol), $ws.Cells.Item($lastRow, $dateCol))
$plainCell = $ws.Cells.Item($firstRow, $dateCol - 2)  # e.g. BD2: untouched, default-styled cell

# Force text formatting first so Excel doesn't reinterpret "2013-06-09" as a
# date serial when it's assigned below.
$dateRange.NumberFormat = "@"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $dateCol)
    if ($cell.Value2 -eq $oldDate) {
        $cell.Value = $newDate
    }
}

# Restore the original (default) cell style now that the text is safely set,
# so the cells keep their original formatting/style index.
$dateRange.Style = $plainCell.Style
